# Natmi following Dr Hou advice
# Update computed NATMI LR-pair statistics for Ntn1-Unc5b sheet (rows 2-10, columns E-T)
# reflecting the ligand/receptor-expressing cell counts moving from 1 to 3 and the
# resulting recomputation of expression / specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$data = @{
    2  = @(3, 1, 2.082653666666667, 6.247961, 0.0472190032704503, 0.0472190032704503, 3, 1, 10.03755233333333, 30.112657, 0.5678353343629902, 0.5678353343629902, 20.90474517137522, 188.142706542377, 0.02681261851036328, 0.02681261851036327)
    3  = @(3, 1, 2.082653666666667, 6.247961, 0.0472190032704503, 0.0472190032704503, 3, 1, 5.808562999999999, 17.425689, 0.3285967737692652, 0.3285967737692652, 12.09722503001433, 108.875025270129, 0.01551601213527035, 0.01551601213527035)
    4  = @(3, 1, 2.082653666666667, 6.247961, 0.0472190032704503, 0.0472190032704503, 3, 1, 1.830756333333333, 5.492269, 0.1035678918677447, 0.1035678918677447, 3.812831390389889, 34.315482513509, 0.004890372624816682, 0.00489037262481668)
    5  = @(3, 1, 37.28222, 111.84666, 0.8452818134314446, 0.8452818134314446, 3, 1, 10.03755233333333, 30.112657, 0.5678353343629902, 0.5678353343629902, 374.2222343528466, 3368.000109175619, 0.479980881160799, 0.479980881160799)
    6  = @(3, 1, 37.28222, 111.84666, 0.8452818134314446, 0.8452818134314446, 3, 1, 5.808562999999999, 17.425689, 0.3285967737692652, 0.3285967737692652, 216.55612364986, 1949.00511284874, 0.2777568768194066, 0.2777568768194066)
    7  = @(3, 1, 37.28222, 111.84666, 0.8452818134314446, 0.8452818134314446, 3, 1, 1.830756333333333, 5.492269, 0.1035678918677447, 0.1035678918677447, 68.25466038572667, 614.29194347154, 0.08754405545123901, 0.08754405545123901)
    8  = @(3, 1, 4.741387, 14.224161, 0.107499183298105, 0.107499183298105, 3, 1, 10.03755233333333, 30.112657, 0.5678353343629902, 0.5678353343629902, 47.59192014508633, 428.327281305777, 0.06104183469182784, 0.06104183469182783)
    9  = @(3, 1, 4.741387, 14.224161, 0.107499183298105, 0.107499183298105, 3, 1, 5.808562999999999, 17.425689, 0.3285967737692652, 0.3285967737692652, 27.540645096881, 247.865805871929, 0.0353238848145882, 0.03532388481458818)
    10 = @(3, 1, 4.741387, 14.224161, 0.107499183298105, 0.107499183298105, 3, 1, 1.830756333333333, 5.492269, 0.1035678918677447, 0.1035678918677447, 8.680324279034334, 78.12291851130901, 0.01113346379168901, 0.01113346379168901)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    # Columns E..T correspond to spreadsheet column indices 5..20
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}
